$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.839.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.37%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.596.80'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -2.55%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''553.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +2.96%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''142.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.31%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.14%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +4.94%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''6.78'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.17%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.101'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -1.96%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  +4.98%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -0.99%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''3.053.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.98%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''58.797.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -1.28%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''20.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -2.28%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''2.581.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -3.86%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E18").Value = '''  -0.47%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''337.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.82%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''10.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.56%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  -0.85%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.10%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''66.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.33%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.428'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +2.66%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.21%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -3.31%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''7.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -2.20%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''0.0₃0753'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +0.87%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -0.08%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''1.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +1.50%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''6.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +2.27%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''154.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +2.28%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''18.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +0.20%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''3.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -2.15%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.894'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +6.90%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -0.94%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''36.84'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = '''0.848'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +1.15%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +1.52%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''3.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.31%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''283.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -1.13%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.996'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -0.37%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.597'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -1.63%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.0956'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +1.03%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D46").Value = '''0.0533'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.38%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -0.25%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''1.941.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -1.33%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''117.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +5.68%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''17.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -2.69%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''4.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -3.18%  '
$ws.Range("E51").Style = "Normal"
